$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 23:05"
$ws.Range("B4").Value = 1615088
$ws.Range("C4").Value = 22365
$ws.Range("D4").Value = 380592
$ws.Range("E4").Value = 1138405
$ws.Range("G4").Value = 1155
$ws.Range("H4").Value = 96091
$ws.Range("B11").Value = 178965
$ws.Range("C11").Value = 434
$ws.Range("E11").Value = 12661
$ws.Range("G11").Value = 34
$ws.Range("H11").Value = 8304
$ws.Range("D28").Value = 27900
$ws.Range("E28").Value = 896
$ws.Range("A93").Value = "Somalia"
$ws.Range("B93").Value = 1594
$ws.Range("C93").Value = 21
$ws.Range("D93").Value = 204
$ws.Range("E93").Value = 1329
$ws.Range("G93").Value = 0
$ws.Range("A94").Value = "Lituania"
$ws.Range("B94").Value = 1593
$ws.Range("C94").Value = 16
$ws.Range("D94").Value = 1049
$ws.Range("E94").Value = 483
$ws.Range("G94").Value = 1
$ws.Range("B128").Value = 585
$ws.Range("C128").Value = 15
$ws.Range("E128").Value = 345
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 35
$ws.Range("A158").Value = "Uganda"
$ws.Range("B158").Value = 160
$ws.Range("D158").Value = 66
$ws.Range("E158").Value = 94
$ws.Range("H158").Value = 0
$ws.Range("A159").Value = "Guadalupe"
$ws.Range("B159").Value = 155
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 109
$ws.Range("E159").Value = 33
$ws.Range("H159").Value = 13
$ws.Range("A160").Value = "Gibraltar"
$ws.Range("B160").Value = 151
$ws.Range("C160").Value = 2
$ws.Range("D160").Value = 146
$ws.Range("E160").Value = 5
